# Auto-generated PowerShell Excel COM-interop script
# Updates the cryptos price list to match the target commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Range, [string]$Text)
    # Force literal text (avoid Excel auto-converting numeric-looking
    # strings like "1.001" or "132.00" into numbers), while keeping
    # the cell's style/formatting identical to before (no text format applied).
    $Range.Value = "'" + $Text
    $Range.Style = "Normal"
}

Set-TextCell $ws.Range("D2") '23.191.11'
Set-TextCell $ws.Range("E2") '  +0.28%  '
Set-TextCell $ws.Range("D3") '1.601.50'
Set-TextCell $ws.Range("E3") '  -0.27%  '
Set-TextCell $ws.Range("D4") '1.001'
Set-TextCell $ws.Range("E4") '  +0.01%  '
Set-TextCell $ws.Range("E5") '  +0.07%  '
Set-TextCell $ws.Range("D6") '303.04'
Set-TextCell $ws.Range("E6") '  +0.24%  '
Set-TextCell $ws.Range("D7") '0.3780'
Set-TextCell $ws.Range("E7") '  -0.34%  '
Set-TextCell $ws.Range("D8") '51.93'
Set-TextCell $ws.Range("E8") '  +3.18%  '
Set-TextCell $ws.Range("D9") '0.3616'
Set-TextCell $ws.Range("E9") '  -1.41%  '
Set-TextCell $ws.Range("D10") '1.266'
Set-TextCell $ws.Range("E10") '  -0.83%  '
Set-TextCell $ws.Range("E11") '  +0.00%  '
Set-TextCell $ws.Range("D12") '0.08114'
Set-TextCell $ws.Range("E12") '  -0.58%  '
Set-TextCell $ws.Range("D13") '22.78'
Set-TextCell $ws.Range("E13") '  -0.64%  '
Set-TextCell $ws.Range("D14") '6.581'
Set-TextCell $ws.Range("E14") '  -0.81%  '
Set-TextCell $ws.Range("D15") '7.399'
Set-TextCell $ws.Range("E15") '  -0.32%  '
Set-TextCell $ws.Range("D16") '0.00001243'
Set-TextCell $ws.Range("E16") '  -1.61%  '
Set-TextCell $ws.Range("D17") '1.599.14'
Set-TextCell $ws.Range("E17") '  -0.21%  '
Set-TextCell $ws.Range("D18") '93.86'
Set-TextCell $ws.Range("E18") '  +1.64%  '
Set-TextCell $ws.Range("D19") '0.06890'
Set-TextCell $ws.Range("E19") '  +0.23%  '
Set-TextCell $ws.Range("D20") '18.06'
Set-TextCell $ws.Range("E20") '  -1.58%  '
Set-TextCell $ws.Range("D21") '6.534'
Set-TextCell $ws.Range("E21") '  -1.29%  '
Set-TextCell $ws.Range("E22") '  -0.02%  '
Set-TextCell $ws.Range("D23") '12.96'
Set-TextCell $ws.Range("E23") '  -1.20%  '
Set-TextCell $ws.Range("D24") '23.184.63'
Set-TextCell $ws.Range("E24") '  +0.28%  '
Set-TextCell $ws.Range("E25") '  +1.90%  '
Set-TextCell $ws.Range("D26") '2.989'
Set-TextCell $ws.Range("E26") '  +6.01%  '
Set-TextCell $ws.Range("D27") '21.21'
Set-TextCell $ws.Range("E27") '  +0.03%  '
Set-TextCell $ws.Range("D28") '149.96'
Set-TextCell $ws.Range("E28") '  -0.05%  '
Set-TextCell $ws.Range("D29") '5.243'
Set-TextCell $ws.Range("E29") '  -0.79%  '
Set-TextCell $ws.Range("D30") '133.65'
Set-TextCell $ws.Range("E30") '  -0.54%  '
Set-TextCell $ws.Range("D31") '2.366'
Set-TextCell $ws.Range("E31") '  -0.93%  '
Set-TextCell $ws.Range("D32") '6.739'
Set-TextCell $ws.Range("E32") '  -2.32%  '
Set-TextCell $ws.Range("D33") '1.778.93'
Set-TextCell $ws.Range("E33") '  -0.13%  '
Set-TextCell $ws.Range("D34") '0.9641'
Set-TextCell $ws.Range("E34") '  -0.02%  '
Set-TextCell $ws.Range("D35") '0.07462'
Set-TextCell $ws.Range("E35") '  -3.68%  '
Set-TextCell $ws.Range("E36") '  -2.29%  '
Set-TextCell $ws.Range("D37") '0.02715'
Set-TextCell $ws.Range("E37") '  -0.97%  '
Set-TextCell $ws.Range("D38") '0.2512'
Set-TextCell $ws.Range("E38") '  -1.85%  '
Set-TextCell $ws.Range("D39") '0.08799'
Set-TextCell $ws.Range("E39") '  -1.17%  '
Set-TextCell $ws.Range("D40") '6.060'
Set-TextCell $ws.Range("E40") '  -4.19%  '
Set-TextCell $ws.Range("B41") 'TheSandbox'
Set-TextCell $ws.Range("C41") 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextCell $ws.Range("D41") '0.7095'
Set-TextCell $ws.Range("E41") '  -0.26%  '
Set-TextCell $ws.Range("B42") 'TrustWalletToken'
Set-TextCell $ws.Range("C42") 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextCell $ws.Range("D42") '1.359'
Set-TextCell $ws.Range("E42") '  -0.88%  '
Set-TextCell $ws.Range("D43") '12.47'
Set-TextCell $ws.Range("E43") '  -1.61%  '
Set-TextCell $ws.Range("D44") '15.54'
Set-TextCell $ws.Range("E44") '  +1.21%  '
Set-TextCell $ws.Range("D45") '0.6526'
Set-TextCell $ws.Range("E45") '  -1.91%  '
Set-TextCell $ws.Range("D46") '2.309'
Set-TextCell $ws.Range("E46") '  -1.01%  '
Set-TextCell $ws.Range("D47") '4.019'
Set-TextCell $ws.Range("E47") '  +0.33%  '
Set-TextCell $ws.Range("D48") '132.00'
Set-TextCell $ws.Range("E48") '  -0.68%  '
Set-TextCell $ws.Range("D49") '0.07946'
Set-TextCell $ws.Range("E49") '  -0.11%  '
Set-TextCell $ws.Range("D50") '1.201'
Set-TextCell $ws.Range("E50") '  -4.32%  '
Set-TextCell $ws.Range("D51") '1.203'
Set-TextCell $ws.Range("E51") '  -0.42%  '
